# Applies the "Revert 'Revert 'Revert 'Yes'''" edit to SUS_INTUI_WORD.xlsx.
#
# Net effect of the commit (stripped of pure Excel-build/locale resave noise
# such as revision GUIDs, rupBuild/calcId numbers, default row-height deltas,
# and chart/theme language tags):
#
#   1. WordHelix!E4 held a stray misspelling ("Comprehensiv"); it is fixed to
#      the already-used "Comprehensive".
#   2. WordHelix had a helper word/count list in columns N:O (rows 1-23) that
#      duplicated the tally already present in column L ("Word = count"
#      strings). That helper list is removed, and column O (rows 2-24) is
#      repurposed to hold just the numeric counts pulled from column L, with
#      O25 keeping its SUM(O2:O24) total.
#   3. WordHeadrush had the same kind of helper list, but in columns O:P
#      (rows 1-28). It is removed, and the numeric counts move into column M
#      (rows 2-28), with M30 keeping a SUM(M2:M28) total (replacing the old
#      SUM(P2:P28) in P30).
#   4. TestDataHelix becomes the active/selected sheet (was WordHeadrush).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. WordHelix!E4 typo fix
# ---------------------------------------------------------------------
$wsWordHelix = $wb.Worksheets.Item("WordHelix")
$wsWordHelix.Range("E4").Value = "Comprehensive"

# ---------------------------------------------------------------------
# 2. WordHelix: drop the old N/O helper list, recompute O from L's counts
# ---------------------------------------------------------------------
$wordHelixCounts = @{
    2 = 2;  3 = 3;  4 = 1;  5 = 1;  6 = 5;  7 = 4;  8 = 1;  9 = 2;
    10 = 2; 11 = 3; 12 = 1; 13 = 1; 14 = 1; 15 = 3; 16 = 1; 17 = 1;
    18 = 1; 19 = 1; 20 = 2; 21 = 1; 22 = 1; 23 = 1; 24 = 1
}

# Clear the redundant helper list (column N rows 1-23, and the old O1).
$wsWordHelix.Range("N1:N23").Value = $null
$wsWordHelix.Range("O1").Value = $null

# Column O rows 2-24 become plain counts (no style), row 25 keeps its
# existing SUM(O2:O24) formula - its cached value recomputes to 40.
foreach ($row in $wordHelixCounts.Keys) {
    $wsWordHelix.Cells.Item($row, 15).Value = $wordHelixCounts[$row]
}
$wsWordHelix.Range("O25").Formula = "=SUM(O2:O24)"

$wsWordHelix.Range("N9").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. WordHeadrush: drop the old O/P helper list, recompute M from L's counts
# ---------------------------------------------------------------------
$wsWordHeadrush = $wb.Worksheets.Item("WordHeadrush")
$wordHeadrushCounts = @{
    2 = 2;  3 = 2;  4 = 1;  5 = 3;  6 = 2;  7 = 1;  8 = 2;  9 = 2;
    10 = 3; 11 = 2; 12 = 1; 13 = 1; 14 = 2; 15 = 2; 16 = 1; 17 = 1;
    18 = 1; 19 = 1; 20 = 2; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1;
    26 = 1; 27 = 1; 28 = 1
}

# Clear the redundant helper list (columns O and P, rows 1-28).
$wsWordHeadrush.Range("O1:P28").Value = $null

# Column M rows 2-28 become plain counts, row 30 gets SUM(M2:M28)
# (replacing the old SUM(P2:P28) in column P).
foreach ($row in $wordHeadrushCounts.Keys) {
    $wsWordHeadrush.Cells.Item($row, 13).Value = $wordHeadrushCounts[$row]
}
$wsWordHeadrush.Range("P30").Value = $null
$wsWordHeadrush.Range("M30").Formula = "=SUM(M2:M28)"

$wsWordHeadrush.Range("P8").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. TestDataHelix becomes the active sheet/selection
# ---------------------------------------------------------------------
$wsTestDataHelix = $wb.Worksheets.Item("TestDataHelix")
$wsTestDataHelix.Activate()
$wsTestDataHelix.Range("Q15").Select() | Out-Null
